# Update cryptos list values (price and 1h volume change) and swap
# the PancakeSwap / Hedera rows, per the Oct 11 2023 GitHub Actions run.
#
# Cells whose new price text is a plain number (e.g. "206.23") are first
# switched to a Text number format so Excel keeps them as strings instead
# of silently converting them to numeric values (which would also lose
# trailing zeros, e.g. "0.930" -> 0.93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.731.41"
$ws.Range("E2").Value = "  -2.67%  "
$ws.Range("D3").Value = "1.565.20"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.23"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.84"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0583"
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "1.562.03"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "26.787.12"
$ws.Range("E16").Value = "  -2.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.33"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.37"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "0.0₃0677"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.29"
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.45"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.72"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0464"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.930"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0162"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.525"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.988"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.33"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.19"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "1.701.09"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.59"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "0.0₇0986"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0948"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -0.81%  "
